$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("N O T E S")

$ws.Range("B6").Value = "employee files its own leave, ot, or ob"

$ws.Range("B6").Select()
